$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new rows (11-13) for the "html" file type, mirroring the
# existing txt/xml/json blocks (English/Arabic/French descriptions).
# The shared-string table records new strings in first-write order, and
# the saved workbook needs them as: html, ملف html, html file,
# Fichier html -- so column A (html) is written first for all rows, then
# the Arabic row's description/lang (B12/C12) is written ahead of the
# English row's (B11/C11), followed by the French row (B13/C13).
$ws.Cells.Item(11, 1).Value = "html"
$ws.Cells.Item(12, 1).Value = "html"
$ws.Cells.Item(13, 1).Value = "html"

$ws.Cells.Item(12, 2).Value = "ملف html"
$ws.Cells.Item(12, 3).Value = "ara"

$ws.Cells.Item(11, 2).Value = "html file"
$ws.Cells.Item(11, 3).Value = "eng"

$ws.Cells.Item(13, 2).Value = "Fichier html"
$ws.Cells.Item(13, 3).Value = "fra"

foreach ($n in 11..13) {
    $ws.Cells.Item($n, 4).Value = $true
    $ws.Cells.Item($n, 4).HorizontalAlignment = -4131
    $ws.Cells.Item($n, 5).Value = "superadmin"
    $ws.Cells.Item($n, 6).Value = "now()"
}

# Update the selection as recorded in the saved file.
$ws.Range("G1:XFD1048576").Select()
